# Adds "Within 5 miles" (col F) and "Within 10 miles" (col G) of HFC
# production facility data to both the "Means" and "Standard Deviations"
# sheets, and updates the "Total Cancer Risk" / "Total Respiratory" rows
# (rows 9 and 10) with refreshed values on both sheets.

$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD    = $wb.Worksheets.Item("Standard Deviations")

# ---------------------------------------------------------------------
# "Means" sheet
# ---------------------------------------------------------------------

# New header cells for the 5-mile / 10-mile columns
$wsMeans.Range("F1").Value = "Within 5 miles of HFC production facility"
$wsMeans.Range("G1").Value = "Within 10 miles of HFC production facility"

# New data columns F (5 mi) and G (10 mi) for existing rows 2-10
$wsMeans.Range("F2").Value = 86
$wsMeans.Range("G2").Value = 79

$wsMeans.Range("F3").Value = 7.4
$wsMeans.Range("G3").Value = 12

$wsMeans.Range("F4").Value = 6.4
$wsMeans.Range("G4").Value = 9.6

$wsMeans.Range("F5").Value = 16
$wsMeans.Range("G5").Value = 19

$wsMeans.Range("F6").Value = 93
$wsMeans.Range("G6").Value = 81

$wsMeans.Range("F7").Value = 3.1
$wsMeans.Range("G7").Value = 4.7

$wsMeans.Range("F8").Value = 2.6
$wsMeans.Range("G8").Value = 3.7

$wsMeans.Range("F9").Value = 26
$wsMeans.Range("G9").Value = 27

$wsMeans.Range("F10").Value = 0.32
$wsMeans.Range("G10").Value = 0.34

# Updated pre-existing values in rows 9 and 10 (Total Cancer Risk / Total
# Respiratory) on the Means sheet
$wsMeans.Range("B9").Value = 29
$wsMeans.Range("C9").Value = 29
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 27

$wsMeans.Range("B10").Value = 0.37
$wsMeans.Range("C10").Value = 0.38
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.31

# ---------------------------------------------------------------------
# "Standard Deviations" sheet
# ---------------------------------------------------------------------

# New header cells for the 5-mile / 10-mile SD columns
$wsSD.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$wsSD.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# New data columns F (5 mi SD) and G (10 mi SD) for existing rows 2-10
$wsSD.Range("F2").Value = 13
$wsSD.Range("G2").Value = 22

$wsSD.Range("F3").Value = 11
$wsSD.Range("G3").Value = 18

$wsSD.Range("F4").Value = 4.4
$wsSD.Range("G4").Value = 8

$wsSD.Range("F5").Value = 9.1
$wsSD.Range("G5").Value = 14

$wsSD.Range("F6").Value = 24
$wsSD.Range("G6").Value = 27

$wsSD.Range("F7").Value = 4.6
$wsSD.Range("G7").Value = 6.2

$wsSD.Range("F8").Value = 4.6
$wsSD.Range("G8").Value = 6.3

$wsSD.Range("F9").Value = 4.6
$wsSD.Range("G9").Value = 4.4

$wsSD.Range("F10").Value = 0.048
$wsSD.Range("G10").Value = 0.05

# Updated pre-existing values in rows 9 and 10 (Total Cancer Risk / Total
# Respiratory) on the Standard Deviations sheet
$wsSD.Range("B9").Value = 10
$wsSD.Range("C9").Value = 7.6
$wsSD.Range("E9").Value = 4.6

$wsSD.Range("B10").Value = 0.14
$wsSD.Range("C10").Value = 0.098
$wsSD.Range("E10").Value = 0.04
